# Update report: a new customer ("Chị duyên") was added to the top of the
# list at LONG XUYÊN branch, pushing the existing rows down by one, and the
# phone/tích lũy/dư nợ information for "Cô tú" was filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 2 (row 1 is the header), shifting all existing
#    customer rows down by one (row 2 -> row 3, ..., row 56 -> row 57).
$ws.Rows.Item(2).Insert()

# 2) Fill in the data for the newly inserted customer row.
$ws.Range("A2").Value = "KH"
$ws.Range("B2").Value = 397
$ws.Range("C2").Value = "Chị duyên"
$ws.Range("D2").Value = "LONG XUYÊN"
$ws.Range("E2").Value = $null
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0988925422"
$ws.Range("G2").Value = $null
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = 8000000
$ws.Range("J2").Value = 0

# 3) Inserting the row turned the previously-empty cells (CCCD/Facebook/Địa
#    chỉ columns, always blank) into zeros for every shifted row. Restore
#    them back to empty, matching the original data shape.
for ($r = 3; $r -le 57; $r++) {
    $ws.Range("E$r").Value = $null
    $ws.Range("G$r").Value = $null
    $ws.Range("H$r").Value = $null
}

# 4) Likewise restore the SĐT (phone) column back to empty for every shifted
#    row that did not originally have a phone number.
$emptyPhoneRows = @(4, 7, 8, 16, 33, 34, 35, 36, 37, 38, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57)
foreach ($r in $emptyPhoneRows) {
    $ws.Range("F$r").Value = $null
}

# 5) "Cô tú" (now on row 39 after the shift) gets her phone number and
#    updated Tích lũy / Dư nợ amounts filled in.
$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = "0939287844"
$ws.Range("I39").Value = 5500000
$ws.Range("J39").Value = 3000000
